$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values stay as text (matching original inline-string cells)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "63.030.38"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "3.040.53"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "588.34"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "151.81"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.538"
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("D9").Value = "3.036.50"
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("D11").Value = "5.79"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").Value = "0.448"
$ws.Range("E12").Value = "  -3.09%  "
$ws.Range("D13").Value = "0.0000235"
$ws.Range("E13").Value = "  -3.01%  "
$ws.Range("D14").Value = "36.32"
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "3.534.93"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").Value = "7.14"
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("D18").Value = "62.964.27"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "3.033.18"
$ws.Range("E19").Value = "  -1.46%  "
$ws.Range("D20").Value = "479.02"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Value = "14.25"
$ws.Range("E21").Value = "  -3.21%  "
$ws.Range("D22").Value = "0.705"
$ws.Range("E22").Value = "  -2.08%  "
$ws.Range("D23").Value = "7.51"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("D24").Value = "2.42"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").Value = "81.93"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").Value = "12.69"
$ws.Range("E26").Value = "  -3.15%  "
$ws.Range("D27").Value = "10.72"
$ws.Range("E27").Value = "  +7.39%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").Value = "7.38"
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "2.67"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "2.20"
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("D33").Value = "27.61"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").Value = "0.0₃0815"
$ws.Range("E36").Value = "  -4.70%  "
$ws.Range("D37").Value = "3.26"
$ws.Range("E37").Value = "  -3.90%  "
$ws.Range("D38").Value = "5.91"
$ws.Range("E38").Value = "  -3.81%  "
$ws.Range("D39").Value = "2.22"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "9.25"
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("D41").Value = "50.40"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "433.70"
$ws.Range("E42").Value = "  -3.16%  "
$ws.Range("D43").Value = "0.287"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("E44").Value = "  +1.85%  "
$ws.Range("D45").Value = "0.0361"
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("D46").Value = "2.819.99"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").Value = "38.31"
$ws.Range("E47").Value = "  -5.94%  "
$ws.Range("D48").Value = "129.10"
$ws.Range("E48").Value = "  -1.61%  "
$ws.Range("D50").Value = "25.10"
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("E51").Value = "  -2.52%  "

# Reset style so no stray number-format style is left attached to cells
$ws.Range("D2:E51").Style = "Normal"

